# Generate Report for Handoff
# Updates the localization-status report with a new handoff run:
#  - new GUID-named source markdown file (was d7fdf396-..., now 9f05bddb-...)
#  - new handoff package hashes for zh-cn / de-de xlf files
#  - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$oldGuid = "d7fdf396-7121-4b4d-9f14-abb869f0b435"
$newGuid = "9f05bddb-1b7c-4bec-bad9-eacb00c554e2"

$newMdName    = "$newGuid.md"
$newZhCnName  = "$newGuid.df656bc531dee5a8c8074fc2c371b51ce6196109.zh-cn.xlf"
$newDeDeName  = "$newGuid.df656bc531dee5a8c8074fc2c371b51ce6196109.de-de.xlf"

# Hyperlink targets stay exactly the same as before the edit - only the
# cell text / hyperlink display text changes.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/6eaa89e2449c0123416d3de9fe19cbae11e8adec/e2e/$oldGuid.md"
$zhCnAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce68c9fef9a6cd531e12eaaa0a05ad38684cb830/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.0bd3b3fb298f17f186f556e59d869eff2ad2df4a.zh-cn.xlf"
$deDeAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1266b5fd4c2b42b4a60cfd3f1eafec5f9beb5fc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.0bd3b3fb298f17f186f556e59d869eff2ad2df4a.de-de.xlf"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)
$wsOverview.Range("D2").Value = "2016-03-23 03:00:59"

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhCnAddress, "", "", $newZhCnName)
$wsZhCn.Range("E2").Value = "2016-03-23 03:00:55"

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deDeAddress, "", "", $newDeDeName)
$wsDeDe.Range("E2").Value = "2016-03-23 03:00:59"
